$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ranked list of GFC inputs (rank 1 through 10)
$values = @(
    "Knowledge Exchange Participation",
    "Grantee-Led Convening Participation",
    "Leveraging",
    "Monitoring Site Visit",
    "Site Visit",
    "Meeting that is not a site visit or at a KE",
    "Phone Call",
    "E-mail",
    "Legal Referral",
    "Additional Touch"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Range("B3").Select()
